$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "equal priors" Laplace/Bayes smoothing matrix (K28:S33 block) used a
# pseudo-count of 2 in L28; the author reworked the sheet to also show the
# "unequal priors" case, and while doing so trimmed the equal-priors
# pseudo-count down to 1. Every L..S formula in rows 30-33/40-44/48-52
# recalculates automatically from this single input cell.
$ws.Range("L28").Value = 1

# Reflect the new scroll position / active cell once the extra matrices
# were added below the original ones.
$ws.Range("L29").Select()

$excel.Calculate()
